$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows before row 16 (everything below shifts down by 8:
# old row 16 -> new row 24, ..., old row 79 -> new row 87)
$ws.Rows("16:23").Insert()

# --- Row 16: Roast year (yyyy) ---
$ws.Cells.Item(16, 1).Value = "~yyyy"
$ws.Cells.Item(16, 2).Value = "Roast year in format yyyy"
$ws.Cells.Item(16, 3).Value = 2020

# --- Row 17: Roast year (yy) ---
$ws.Cells.Item(17, 1).Value = "~yy"
$ws.Cells.Item(17, 2).Value = "Roast year in format yy"
$ws.Cells.Item(17, 3).Value = 20

# --- Row 18: Roast month (MMM, localized) ---
$ws.Cells.Item(18, 1).Value = "~mmm"
$ws.Cells.Item(18, 3).Value = "Feb"
$ws.Cells.Item(18, 2).Value = "Roast month in format MMM (localized)"

# --- Row 19: Roast month (MM) ---
$ws.Cells.Item(19, 2).Value = "Roast month in format MM"
$ws.Cells.Item(19, 3).Value = "02"
$ws.Cells.Item(19, 1).Value = "~mm"

# --- Row 20: Roast day (ddd, localized) ---
$ws.Cells.Item(20, 1).Value = "~ddd"
$ws.Cells.Item(18, 2).Value = "Roast month in format MMM (localized)"
$ws.Cells.Item(20, 2).Value = "Roast day in format ddd (localized)"
$ws.Cells.Item(20, 3).Value = "Wed"

# --- Row 21: Roast day (dd) ---
$ws.Cells.Item(21, 2).Value = "Roast day in format dd"
$ws.Cells.Item(21, 1).Value = "~dd"
$ws.Cells.Item(21, 3).Value = "05"

# --- Row 22: Roast hour (hh) ---
$ws.Cells.Item(22, 1).Value = "~hour"
$ws.Cells.Item(23, 1).Value = "~minute"
$ws.Cells.Item(22, 2).Value = "Roast hour in format hh"
$ws.Cells.Item(23, 2).Value = "Roast minute in format mm"

# --- Row 22/23 numeric examples ---
$ws.Cells.Item(22, 3).Value = 17
$ws.Cells.Item(23, 3).Value = 42

# Apply the same "quote-prefixed text" style (style index 5 in styles.xml)
# used by other localized/text example cells to C19 and C21.
$ws.Range("C19").Style = $ws.Range("C40").Style
$ws.Range("C21").Style = $ws.Range("C40").Style

# Update view: scroll so row 5 is near top, and select C24 (what used to be
# row 16, now shifted to row 24) as the active cell.
$ws.Range("C24").Select()
